$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Change 1: Collapse the three long "CORE COMPETENCIES" paragraphs
# into a single summary paragraph.
# ------------------------------------------------------------------
$bullet = [char]0x2022

$coreCompIndex = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "Data Visualization & Design:*") {
        $coreCompIndex = $i
        break
    }
}

if ($coreCompIndex -ne $null) {
    $firstPara = $d.Paragraphs.Item($coreCompIndex)
    $firstPara.Range.Text = "Data Visualization & Design " + $bullet + " Geospatial Analysis & Mapping " + $bullet + " Technical Visualization"

    # The two following paragraphs (previously "Geospatial Analysis & Mapping: ..."
    # and "Technical Visualization: ...") are no longer needed.
    $d.Paragraphs.Item($coreCompIndex + 1).Range.Delete()
    $d.Paragraphs.Item($coreCompIndex + 1).Range.Delete()
}

# ------------------------------------------------------------------
# Change 2: Add a new "TECHNICAL SKILLS" section right before the
# closing "For a more detailed..." paragraph.
# ------------------------------------------------------------------
$ledIndex = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "*Led multi-million dollar*") {
        $ledIndex = $i
        break
    }
}

if ($ledIndex -ne $null) {
    $ledPara = $d.Paragraphs.Item($ledIndex)
    $ledPara.Range.InsertParagraphAfter()

    $heading = $d.Paragraphs.Item($ledIndex + 1)
    $heading.Range.Text = "TECHNICAL SKILLS"
    $heading.Style = "Heading2"

    $heading.Range.InsertParagraphAfter()
    $line1 = $d.Paragraphs.Item($ledIndex + 2)
    $line1.Style = "Normal"
    $line1.Range.Text = "DATA VISUALIZATION & DESIGN Interactive Dashboards; Statistical Visualization; Geospatial Mapping; Choropleth Design"

    $line1.Range.InsertParagraphAfter()
    $line2 = $d.Paragraphs.Item($ledIndex + 3)
    $line2.Style = "Normal"
    $line2.Range.Text = "GEOSPATIAL ANALYSIS & MAPPING Spatial Analysis; Mapping Technologies; Web Mapping; Spatial Data Processing"

    $line2.Range.InsertParagraphAfter()
    $line3 = $d.Paragraphs.Item($ledIndex + 4)
    $line3.Style = "Normal"
    $line3.Range.Text = "TECHNICAL VISUALIZATION Programming; Database Integration; Web Technologies; Statistical Computing"
}
